$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains one new data row: a new record is inserted at row 117,
# pushing the former rows 117-327 down to 118-328 (dimension grows to R328).
# The new row 117 duplicates the data of the (now shifted) row below it,
# except for its date (column D), which is a new value (44725 = 2022-06-13).

$ws.Rows("117:117").Insert()

$ws.Range("A117").Value = $ws.Range("A118").Value2
$ws.Range("B117").Value = $ws.Range("B118").Value2
$ws.Range("C117").Value = $ws.Range("C118").Value2
$ws.Range("D117").Value = 44725
$ws.Range("E117").Value = $ws.Range("E118").Value2
$ws.Range("F117").Value = $ws.Range("F118").Value2
$ws.Range("G117").Value = $ws.Range("G118").Value2
$ws.Range("H117").Value = $ws.Range("H118").Value2
$ws.Range("I117").Value = $ws.Range("I118").Value2
$ws.Range("J117").Value = $ws.Range("J118").Value2
$ws.Range("K117").Value = $ws.Range("K118").Value2
$ws.Range("L117").Value = $ws.Range("L118").Value2
$ws.Range("M117").Value = $ws.Range("M118").Value2
$ws.Range("N117").Value = $ws.Range("N118").Value2
$ws.Range("O117").Value = $ws.Range("O118").Value2
$ws.Range("P117").Value = $ws.Range("P118").Value2
$ws.Range("Q117").Value = $ws.Range("Q118").Value2
$ws.Range("R117").Value = $ws.Range("R118").Value2
